$wb = $excel.ActiveWorkbook

# Sheet ALC, row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 9999
$ws.Range("I47").Value = 9999
$ws.Range("K47").Value = 9999
$ws.Range("M47").Value = -9027

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4242.5
$ws.Range("J62").Value = 3494
$ws.Range("L62").Value = 3494
$ws.Range("N62").Value = -4742

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4242.5
$ws.Range("J65").Value = 3494
$ws.Range("L65").Value = 17470
$ws.Range("N65").Value = -23710

# Sheet ALC, row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1102.4
$ws.Range("I99").Value = 1332.75
$ws.Range("J99").Value = 181
$ws.Range("K99").Value = 3998.25
$ws.Range("L99").Value = 543
$ws.Range("M99").Value = -2500.25
$ws.Range("N99").Value = -3539

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10724.75
$ws.Range("I61").Value = 6949.5
$ws.Range("K61").Value = 6949.5
$ws.Range("M61").Value = -6737.5

# Sheet ARM, row 124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 59497.75
$ws.Range("J124").Value = 59497.75
$ws.Range("L124").Value = 59497.75
$ws.Range("N124").Value = -69317.75

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 18855.572
$ws.Range("I132").Value = 17336.666
$ws.Range("K132").Value = 52009.99800000001
$ws.Range("M132").Value = -49479.99800000001

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10724.75
$ws.Range("I136").Value = 6949.5
$ws.Range("K136").Value = 20848.5
$ws.Range("M136").Value = -18298.5

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1880.3334
$ws.Range("I134").Value = 1876.6
$ws.Range("J134").Value = 1899
$ws.Range("K134").Value = 5629.799999999999
$ws.Range("L134").Value = 5697
$ws.Range("M134").Value = -3094.799999999999
$ws.Range("N134").Value = -10767

# Sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 160
$ws.Range("I2").Value = 187.04546
$ws.Range("J2").Value = 11.25
$ws.Range("K2").Value = 1122.27276
$ws.Range("L2").Value = 67.5
$ws.Range("M2").Value = -1009.27276
$ws.Range("N2").Value = -293.5

# Sheet CUL, row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 100.25
$ws.Range("I7").Value = 100.5
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 301.5
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -189.5
$ws.Range("N7").Value = -524

# Sheet CUL, row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 97.14286
$ws.Range("I23").Value = 67.5
$ws.Range("K23").Value = 202.5
$ws.Range("M23").Value = 32.5

# Sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 963.63635
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 963.63635
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 2890.90905
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -3058.90905

# Sheet CUL, row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 60.6
$ws.Range("I38").Value = 73.333336
$ws.Range("J38").Value = 41.5
$ws.Range("K38").Value = 220.000008
$ws.Range("L38").Value = 124.5
$ws.Range("M38").Value = 126.999992
$ws.Range("N38").Value = -818.5

# Sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 100
$ws.Range("I39").Value = 100
$ws.Range("K39").Value = 300
$ws.Range("M39").Value = -6

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = $null

# Sheet CUL, row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 611.875
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 798.3333
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 2394.9999
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6234.9999

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1016.2
$ws.Range("I140").Value = 1016.2
$ws.Range("K140").Value = 3048.6
$ws.Range("M140").Value = 2131.4

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15669.667
$ws.Range("I70").Value = 12000
$ws.Range("K70").Value = 12000
$ws.Range("M70").Value = -11730

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 15669.667
$ws.Range("I73").Value = 12000
$ws.Range("K73").Value = 12000
$ws.Range("M73").Value = -11064

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3499
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3499
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3499
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -5495

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3499
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3499
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17495
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -27479

# Sheet GSM, row 99
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 10000
$ws.Range("K99").Value = 10000
$ws.Range("M99").Value = -7754

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 412399.6
$ws.Range("I132").Value = 511249.5
$ws.Range("K132").Value = 1533748.5
$ws.Range("M132").Value = -1531218.5

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6500
$ws.Range("J7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("N7").Value = -7224

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1200
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -1472

# Sheet LTW, row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 25000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 25000
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = -26996

# Sheet LTW, row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 25000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = -84984

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6500
$ws.Range("J126").Value = 7000
$ws.Range("L126").Value = 21000
$ws.Range("N126").Value = -25940

# Sheet LTW, row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11167.333
$ws.Range("I132").Value = 6334.6665
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 19003.9995
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -16473.9995
$ws.Range("N132").Value = -53060

# Sheet WVR, row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1699.8
$ws.Range("I100").Value = 1399.5
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 2799
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -2258
$ws.Range("N100").Value = -4882
